$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix job "text color" (column I) for a few rows: white -> black ---
$ws.Range("I2").Value = "black"
$ws.Range("I11").Value = "black"
$ws.Range("I12").Value = "black"
$ws.Range("I15").Value = "black"

# --- Fix job descriptions (column J) ---

# Row 3: Warner Brothers Interactive Entertainment - fix "(ZenDesk](url)" typo -> "[ZenDesk](url)"
$ws.Range("J3").Value = '• Worked as a consulting data engineer in the analytics team of [Warner Brothers Interactive Entertainment](url) Division.
• Implemented high-volume pipeline integrations shuffling game telemetry and [user PII data](URL) between WB-distributed consumer games and [marketing service platforms](URL) via [Segment customer data platform](url) using [Kafka](url), [Redshift](url), and [Airflow](url). 
• Employed [Python](url), [Amazon Glue](url) and [Apache Airflow](url) for external [3rd-party integrations](url) and internal [dev-ops integrations](url) with [Jenkins](url), [DataDog](url), and [ZenDesk](url) 
• Integrated with [Google BigQuery](url) [data warehouse](url) and [AWS-managed services](url) [Airflow](url), [S3](url), [Glue](url), and [Redshift warehouse](url).'

# Row 4: Angel Studios - add missing space "network](url)with" -> "network](url) with"
$ws.Range("J4").Value = '• Worked as a consulting data engineer for [Angel Studios](url], a streaming media service that offers family-friendly entertainment that amplifies light, with titles including The Chosen, Dry Bar Comedy, and Tuttle Twins.
• Used [Python](url), [Pandas](url]), [Numpy](url), [Keras](url), and [Jupyter](url) to build and tune [hyperparameters](url) of a [convolutional neural network](url) with [supervised learning](url) on [AWS Sagemaker](url) to classify movie frames from episodic programs stored in [S3](url). 
• Built web client apps using [Python](url) with [Postman](url) that made [RESTful API](url) requests to pull monthly usage data from various web marketing partners like [FaceBook](url), [Google Play](url), and [Vimeo](url). 
• Worked with [Segment customer data platform](url), [Excel](url), and [Tableau](url) to create scheduled [reports] for the company''s sales and finance teams.'

# Row 11: HomePortfolio LLC - fix "and Omnigraph](url)" -> "and [Omnigraph](url)"
$ws.Range("J11").Value = '• As co-founder and CTO, designed and led the development of a public website used by discerning home designers and builders called HomePortfolio.com. 
• Hired a staff of   10 software and databases developers. 
• Worked with data acquisition team to scan and tag over 700,000 premium home design products from over 2,000 manufactures and vendors.
• Designed datamodel and data entry tools for category-specific product attribution. 
• Helped extended the business model to provide online product selection tools for participating vendors and manufacturers.
• Instrumental in raising over $70M in venture capital. 
• Used [Oracle IIi](url), [ATG Dyanamo](url), [Java 8](url), [Akamai](url), [WebTrends](url), [FileMakerPro](url), [ImageMagik](url), and [Omnigraph](url) for data modeling and workflow designs.'
